# "complete monthly and re-run daily"
# Update arrive_thres counts (and their dependent ratios) for two team
# members (sp99063 / 張宜君 and sp99025 / 林川評) across the four sheets
# of the workbook: team_df, team_df_day, productivity_tl and
# productivity_team_function.

$wb = $excel.ActiveWorkbook

# --- team_df (daily detail, one row per person per shift) ---
$wsTeamDf = $wb.Worksheets.Item("team_df")

# Row 6 -> sp99063 / 張宜君
$wsTeamDf.Range("S6").Value = 5
$wsTeamDf.Range("U6").Value = $wsTeamDf.Range("S6").Value / $wsTeamDf.Range("T6").Value

# Row 7 -> sp99025 / 林川評
$wsTeamDf.Range("S7").Value = 4
$wsTeamDf.Range("U7").Value = $wsTeamDf.Range("S7").Value / $wsTeamDf.Range("T7").Value

# --- team_df_day (per person per day aggregate) ---
$wsTeamDfDay = $wb.Worksheets.Item("team_df_day")

# Row 3 -> sp99063 / 張宜君
$wsTeamDfDay.Range("F3").Value = 8
$wsTeamDfDay.Range("H3").Value = $wsTeamDfDay.Range("F3").Value / $wsTeamDfDay.Range("G3").Value

# Row 7 -> sp99025 / 林川評
$wsTeamDfDay.Range("F7").Value = 11
$wsTeamDfDay.Range("H7").Value = $wsTeamDfDay.Range("F7").Value / $wsTeamDfDay.Range("G7").Value

# --- productivity_tl (TL productivity score, mirrors team_df_day ratio) ---
$wsProdTl = $wb.Worksheets.Item("productivity_tl")

$wsProdTl.Range("D3").Value = $wsTeamDfDay.Range("H3").Value
$wsProdTl.Range("D7").Value = $wsTeamDfDay.Range("H7").Value

# --- productivity_team_function (mirrors team_df_day ratio as well) ---
$wsProdFunc = $wb.Worksheets.Item("productivity_team_function")

$wsProdFunc.Range("D3").Value = $wsTeamDfDay.Range("H3").Value
$wsProdFunc.Range("D7").Value = $wsTeamDfDay.Range("H7").Value
